$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new row before row 172 (pushes the old row 172.. down by one)
$ws.Rows.Item(172).Insert()

# The insert carried formatting from the row above into every column,
# including column G which should stay fully empty (no cell at all).
$ws.Cells.Item(172, 7).Clear()

# Populate the new acronym-key row: indst / ItUBB / Industries that Use Byproduct Biomass / low
$ws.Cells.Item(172, 1).Value = "indst"
$ws.Cells.Item(172, 2).Value = "ItUBB"
$ws.Cells.Item(172, 3).Value = "Industries that Use Byproduct Biomass"
$ws.Cells.Item(172, 6).Value = "low"

# Match the "low" fill/style used elsewhere in column F (row 6 is a template with that style)
$ws.Cells.Item(6, 6).Copy()
$ws.Cells.Item(172, 6).PasteSpecial(-4122)

# Update the saved selection to reflect the new row location
$ws.Range("A172").Select()
